$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table row 1): 54÷3=, 57÷9=, 64÷8=, 52÷7=, 95÷8=
$t.Cell(1, 1).Range.Text = "19÷7="
$t.Cell(1, 2).Range.Text = "55÷9="
$t.Cell(1, 3).Range.Text = "27÷4="
$t.Cell(1, 4).Range.Text = "61÷2="
$t.Cell(1, 5).Range.Text = "31÷6="

# Row 2 (table row 5): 66÷8=, 78÷7=, 24÷6=, 87÷6=, 82÷8=
$t.Cell(5, 1).Range.Text = "67÷9="
$t.Cell(5, 2).Range.Text = "77÷2="
$t.Cell(5, 3).Range.Text = "99÷8="
$t.Cell(5, 4).Range.Text = "66÷9="
$t.Cell(5, 5).Range.Text = "34÷3="

# Row 3 (table row 9): 61÷8=, 18÷9=, 18÷3=, 22÷6=, 34÷3=
$t.Cell(9, 1).Range.Text = "64÷9="
$t.Cell(9, 2).Range.Text = "81÷7="
$t.Cell(9, 3).Range.Text = "40÷6="
$t.Cell(9, 4).Range.Text = "28÷5="
$t.Cell(9, 5).Range.Text = "57÷6="

# Row 4 (table row 13): 51÷6=, 86÷4=, 96÷2=, 18÷2=, 10÷8=
$t.Cell(13, 1).Range.Text = "26÷8="
$t.Cell(13, 2).Range.Text = "88÷9="
$t.Cell(13, 3).Range.Text = "80÷6="
$t.Cell(13, 4).Range.Text = "46÷7="
$t.Cell(13, 5).Range.Text = "96÷8="

# Row 5 (table row 17): 21÷3=, 60÷3=, 34÷9=, 92÷6=, 38÷8=
$t.Cell(17, 1).Range.Text = "84÷3="
$t.Cell(17, 2).Range.Text = "11÷4="
$t.Cell(17, 3).Range.Text = "14÷7="
$t.Cell(17, 4).Range.Text = "44÷8="
$t.Cell(17, 5).Range.Text = "90÷7="
